$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasting Notes")

# --- Row 11: new tasting-note entry (mirrors the formatting of rows 8-10) ---

# Clone formatting from row 10 (date/day#/score/notes row) onto row 11 first,
# so the new row picks up the same borders / alignment / wrap / number format.
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Date (05 Dec 2020 -> serial 44170)
$ws.Range("A11").Value = 44170

# Score
$ws.Range("C11").Value = 3

# Notes (new shared string)
$ws.Range("D11").Value = "Cooled and served @ 12 C. Opening the swing-top gave a very small puff. Poured clear with a little foam. Light carbonation. Very light haze due to yeast not fully deposited. Light body, crisp and hoppy and a mild grassy/vegetal taste and a dry mouthfeel."

# Day # formula
$ws.Range("B11").Formula = '=A11-$A$6'

# Re-apply General/centered/bordered formatting to B11 (writing the formula
# above can cause the engine to auto-infer a date format from its operands).
$ws.Range("C11").Copy()
$ws.Range("B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match row height used by the other wrapped-notes rows.
$ws.Rows(11).RowHeight = 60

# --- Selection / scrolled view ---
$ws.Range("A12").Select()
